# Update the "Industries" column (H) for rows 22 through 91 from 1 to 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H22:H91").Value = 0
